$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IAR-DST")
$ws.Range("O2:O11").ClearContents()
